$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: "001" -> "003" (keep as text so leading zero is preserved)
$ws.Range("J2").Value = "'003"

# N2: report date changes from 2019-12-31 to 2020-03-31 (stored as text, like original)
$ws.Range("N2").Value = "2020-03-31 00:00:00"

# Numeric metric updates for row 2
$ws.Range("O2").Value = 410193440.08
$ws.Range("P2").Value = 80583773.5
$ws.Range("Q2").Value = 176966882.19
$ws.Range("R2").Value = 30.4781761532
$ws.Range("S2").Value = 33030334.57
$ws.Range("T2").Value = -13.2021228327
$ws.Range("U2").Value = 46751824.1
$ws.Range("V2").Value = -6.6350670126
$ws.Range("W2").Value = 38280502.37
$ws.Range("X2").Value = 20025289.59
$ws.Range("Y2").Value = -6.4164159927

# Z2 and AA2 no longer carry values - clear them to blank
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()

$ws.Range("AB2").Value = 371912937.71
$ws.Range("AC2").Value = 13.6067368189
$ws.Range("AD2").Value = 8.363779900200001
$ws.Range("AE2").Value = -25.1822166766
$ws.Range("AF2").Value = 891.2122964855
$ws.Range("AG2").Value = 9.332304866299999
